# Apply crypto price/volume updates (and a row49/row50 coin swap)
# as captured by the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value2 = '27.664.48'
$ws.Cells.Item(2, 5).Value2 = '  -0.46%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value2 = '1.848.42'
$ws.Cells.Item(3, 5).Value2 = '  -0.82%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value2 = '1.013'
$ws.Cells.Item(4, 5).Value2 = '  -1.85%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '318.72'
$ws.Cells.Item(5, 5).Value2 = '  -1.31%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '1.009'
$ws.Cells.Item(6, 5).Value2 = '  -1.87%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '0.4293'
$ws.Cells.Item(7, 5).Value2 = '  -2.49%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = '0.3738'
$ws.Cells.Item(8, 5).Value2 = '  -1.95%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = '0.07327'
$ws.Cells.Item(9, 5).Value2 = '  -1.67%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '0.8753'
$ws.Cells.Item(10, 5).Value2 = '  -1.52%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '21.55'
$ws.Cells.Item(11, 5).Value2 = '  -0.66%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = '1.847.81'
$ws.Cells.Item(12, 5).Value2 = '  -1.77%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = '6.715'
$ws.Cells.Item(13, 5).Value2 = '  -0.75%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value2 = '5.425'
$ws.Cells.Item(14, 5).Value2 = '  -2.33%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = '0.07101'

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '88.92'
$ws.Cells.Item(16, 5).Value2 = '  +4.07%  '

# Row 17
$ws.Cells.Item(17, 5).Value2 = '  -2.03%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = '0.000008965'
$ws.Cells.Item(18, 5).Value2 = '  -1.90%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = '1.010'
$ws.Cells.Item(19, 5).Value2 = '  -1.76%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '15.42'
$ws.Cells.Item(20, 5).Value2 = '  -0.98%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = '27.685.85'
$ws.Cells.Item(21, 5).Value2 = '  -0.45%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '5.202'
$ws.Cells.Item(22, 5).Value2 = '  -2.10%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = '11.04'
$ws.Cells.Item(23, 5).Value2 = '  -2.17%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = '2.076.55'
$ws.Cells.Item(24, 5).Value2 = '  -1.35%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '1.988'
$ws.Cells.Item(25, 5).Value2 = '  -1.90%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = '154.95'
$ws.Cells.Item(26, 5).Value2 = '  -2.31%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '18.60'
$ws.Cells.Item(27, 5).Value2 = '  -1.08%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '2.186'
$ws.Cells.Item(28, 5).Value2 = '  +9.37%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = '5.363'
$ws.Cells.Item(29, 5).Value2 = '  -0.68%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = '118.80'
$ws.Cells.Item(30, 5).Value2 = '  +0.65%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = '0.08925'
$ws.Cells.Item(31, 5).Value2 = '  -1.11%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '1.225'
$ws.Cells.Item(32, 5).Value2 = '  +0.13%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '0.7758'
$ws.Cells.Item(33, 5).Value2 = '  -1.20%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '4.542'
$ws.Cells.Item(34, 5).Value2 = '  -1.07%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '2.924'
$ws.Cells.Item(35, 5).Value2 = '  -3.36%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = '1.010'
$ws.Cells.Item(36, 5).Value2 = '  -1.96%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '1.128'
$ws.Cells.Item(37, 5).Value2 = '  -1.54%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = '0.01978'
$ws.Cells.Item(38, 5).Value2 = '  -0.46%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '0.05333'
$ws.Cells.Item(39, 5).Value2 = '  +0.10%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '7.272'
$ws.Cells.Item(40, 5).Value2 = '  +5.03%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = '2.905'
$ws.Cells.Item(41, 5).Value2 = '  +1.62%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '0.1688'

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = '0.5116'
$ws.Cells.Item(43, 5).Value2 = '  -2.13%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = '8.790'
$ws.Cells.Item(44, 5).Value2 = '  -1.37%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '10.67'
$ws.Cells.Item(45, 5).Value2 = '  -0.51%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '108.02'
$ws.Cells.Item(46, 5).Value2 = '  -2.71%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '0.4771'
$ws.Cells.Item(47, 5).Value2 = '  +0.58%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '0.06462'
$ws.Cells.Item(48, 5).Value2 = '  -2.17%  '

# Row 49
$ws.Cells.Item(49, 2).Value2 = 'PaxDollar'
$ws.Cells.Item(49, 3).Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = '1.010'
$ws.Cells.Item(49, 5).Value2 = '  -2.05%  '

# Row 50
$ws.Cells.Item(50, 2).Value2 = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = '1.686'
$ws.Cells.Item(50, 5).Value2 = '  -2.30%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = '1.835'
$ws.Cells.Item(51, 5).Value2 = '  -4.58%  '
